# Remove the old index column (A) which held plain numeric row indices
# 0..7, shifting B:J left to A:I, then give the new header cell (A1) the
# same header style as its neighbours plus the "Datasets" label so the
# label column keeps its header while the numeric index column is gone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting column A shifts B:J -> A:I, carrying along each cell's value,
# style and type (numbers stay numbers, shared strings stay strings).
$ws.Range("A1").EntireColumn.Delete()

# The new column A (former labels column B) picks up the header style
# used by the rest of row 1, and gets its own header text restored.
$ws.Range("A1").Value = "Datasets"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# The label cells A2:A9 (former B2:B9) should no longer carry the bordered
# header-ish style that the old index column A used to have.
$ws.Range("B2").Copy()
$ws.Range("A2:A9").PasteSpecial(-4122)

$excel.CutCopyMode = 0
